# Updated cryptos list on Fri Jun 28 04:46:58 UTC 2024 with GitHub Actions
#
# The sheet stores Coin/Link/Price/Volume(1h) as plain text cells (many of
# the "Price" values look numeric, e.g. "9.54" or have multiple dots as a
# thousands separator, e.g. "61.624.72"). To keep them stored as text
# (rather than letting Excel auto-convert a numeric-looking string into a
# real number) we force the NumberFormat to "@" (Text) just before writing
# any value that would otherwise parse as a plain number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; D='61.624.72'; E='  +0.78%  '},
    @{Row=3; D='3.451.00'; E='  +1.88%  '},
    @{Row=4; E='  -0.04%  '},
    @{Row=5; D='581.15'; E='  +1.19%  '},
    @{Row=6; D='145.93'; E='  +6.33%  '},
    @{Row=7; D='3.452.28'; E='  +1.98%  '},
    @{Row=8; E='  +0.03%  '},
    @{Row=9; E='  +1.43%  '},
    @{Row=10; E='  +0.07%  '},
    @{Row=11; E='  +2.70%  '},
    @{Row=12; E='  +2.34%  '},
    @{Row=13; D='4.040.13'; E='  +2.20%  '},
    @{Row=14; D='28.05'; E='  +9.20%  '},
    @{Row=15; E='  -0.99%  '},
    @{Row=16; E='  +1.19%  '},
    @{Row=17; D='3.466.74'; E='  +2.42%  '},
    @{Row=18; D='61.739.34'; E='  +0.75%  '},
    @{Row=19; E='  +8.42%  '},
    @{Row=20; E='  +3.66%  '},
    @{Row=21; D='9.54'; E='  +2.18%  '},
    @{Row=22; D='389.99'; E='  +3.52%  '},
    @{Row=23; E='  +2.86%  '},
    @{Row=24; D='73.72'; E='  +3.83%  '},
    @{Row=25; E='  +0.10%  '},
    @{Row=26; D='0.995'; E='  -0.60%  '},
    @{Row=27; D='0.0000124'; E='  -1.35%  '},
    @{Row=28; D='3.595.83'; E='  +2.21%  '},
    @{Row=29; E='  +1.81%  '},
    @{Row=30; E='  +2.69%  '},
    @{Row=31; E='  +0.20%  '},
    @{Row=32; D='8.20'; E='  +1.60%  '},
    @{Row=33; D='1.47'; E='  -10.85%  '},
    @{Row=34; E='  +2.18%  '},
    @{Row=36; E='  +2.80%  '},
    @{Row=37; D='3.478.12'; E='  +2.06%  '},
    @{Row=38; D='7.03'; E='  +2.80%  '},
    @{Row=39; E='  +0.81%  '},
    @{Row=40; E='  -0.25%  '},
    @{Row=41; D='166.96'; E='  +1.34%  '},
    @{Row=42; D='0.0786'; E='  +3.23%  '},
    @{Row=43; D='27.51'; E='  +7.79%  '},
    @{Row=44; E='  +3.93%  '},
    @{Row=45; E='  +3.99%  '},
    @{Row=46; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='1.00'; E='  +0.02%  '},
    @{Row=47; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='42.45'; E='  +1.73%  '},
    @{Row=49; E='  -2.28%  '},
    @{Row=50; D='2.572.90'; E='  +0.58%  '},
    @{Row=51; D='6.97'; E='  +2.57%  '}
)

foreach ($item in $changes) {
    $row = $item.Row

    if ($item.ContainsKey("B")) {
        $ws.Range("B$row").Value = $item.B
    }
    if ($item.ContainsKey("C")) {
        $ws.Range("C$row").Value = $item.C
    }
    if ($item.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        # Keep the cell as text: if the new price string would otherwise be
        # auto-recognised as a plain number, force a Text number format
        # first so it round-trips the same way the original data did.
        if ($item.D -match '^[+-]?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E$row").Value = $item.E
    }
}
